{"js": "// After the final paragraph (\"... to avoid diverging.\"), append:\n//   - three blank paragraphs\n//   - \"I should disable right clicking context menu when it\u2019s over a link\"\n//   - \"Need to make sure flowgraph names are unique within composite\"\nconst body = context.document.body;\n\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\n  \"I should disable right clicking context menu when it\\u2019s over a link\",\n  Word.InsertLocation.end\n);\nbody.insertParagraph(\n  \"Need to make sure flowgraph names are unique within composite\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# After the final paragraph (\"... to avoid diverging.\"), append:\n#   - three blank paragraphs\n#   - \"I should disable right clicking context menu when it's over a link\"\n#   - \"Need to make sure flowgraph names are unique within composite\"\n$d = $word.ActiveDocument\n\n$d.Paragraphs.Add() | Out-Null\n$d.Paragraphs.Add() | Out-Null\n$d.Paragraphs.Add() | Out-Null\n\n$p4 = $d.Paragraphs.Add()\n$p4.Range.Text = \"I should disable right clicking context menu when it\u2019s over a link\"\n\n$p5 = $d.Paragraphs.Add()\n$p5.Range.Text = \"Need to make sure flowgraph names are unique within composite\"\n"}
